$d = $word.ActiveDocument

# --- Paragraph 1: "...thành 25/05/2012." -> "...thành 29/05/2012." ---
# Split into 3 runs: "...thành ", "29", "/05/2012."
$m1 = $d.Content.Duplicate
$null = $m1.Find.Execute("25/05/2012", $false, $false, $false, $false, $false, `
                          $true, 1, $false, "", 0)

$date1 = $d.Range($m1.Start, $m1.End)
$num1 = $d.Range($date1.Start, $date1.Start + 2)
$num1.Text = "29"

# Re-seat after the text mutation, then force a run split at both edges of
# the "29" substring by toggling a character formatting property on/off.
$num1b = $d.Range($date1.Start, $date1.Start + 2)
$num1b.Bold = 1
$num1b.Bold = 0

# --- Paragraph 2: "...đã trở thành  25/05/2012." -> "...đã trở thành  29/05/2012." ---
# Split into 3 runs: "...đã trở", " thành  29", "/05/2012."
$m2 = $d.Content.Duplicate
$m2.Start = $m1.End
$null = $m2.Find.Execute("25/05/2012", $false, $false, $false, $false, $false, `
                          $true, 1, $false, "", 0)

$date2 = $d.Range($m2.Start, $m2.End)
$num2 = $d.Range($date2.Start, $date2.Start + 2)
$num2.Text = "29"

$split2 = $d.Range($date2.Start - 8, $date2.Start + 2)
$split2.Bold = 1
$split2.Bold = 0
